$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.428.08"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "2.162.84"
$ws.Range("E3").Value = "  +3.68%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.20"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.11"
$ws.Range("E7").Value = "  +4.77%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0861"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.03"
$ws.Range("E12").Value = "  +7.07%  "
$ws.Range("D13").Value = "2.481.07"
$ws.Range("E13").Value = "  +3.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.23"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.815"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.57"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "2.160.96"
$ws.Range("E17").Value = "  +3.95%  "
$ws.Range("D18").Value = "39.431.88"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.22"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.13"
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("D21").Value = "0.0₃0853"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.75"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  -2.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.64"
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.16"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.72"
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("E30").Value = "  -3.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.62"
$ws.Range("E31").Value = "  +9.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("E33").Value = "  +3.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.81"
$ws.Range("E34").Value = "  +2.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.09"
$ws.Range("E35").Value = "  +9.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0623"
$ws.Range("E36").Value = "  +2.53%  "
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.59"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "104.36"
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.02"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").Value = "1.537.94"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("E44").Value = "  +6.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0932"
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.89"
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  +7.45%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.81"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("D50").Value = "2.366.05"
$ws.Range("E50").Value = "  +3.59%  "
$ws.Range("E51").Value = "  +0.35%  "

Write-Host "Update complete"
